$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E11").Value = 156884
$ws.Range("C13").Value = 6
$ws.Range("F13").Value = 1

$ws.Rows("22:22").Delete()
